# Apply updated "Zeile im Spreadsheet" (column K) values on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    5  = 505
    6  = 415
    8  = 439
    9  = 328
    10 = 517
    12 = 429
    14 = 510
    17 = 423
    21 = 508
    23 = 443
    24 = 330
    25 = 453
    26 = 318
    27 = 519
}

foreach ($row in $updates.Keys) {
    $ws.Range("K$row").Value = $updates[$row]
}
